# Insert a new "Industry" column before column C ("Mutual Fund"), shifting
# the existing Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/MoM/QoQ columns
# one place to the right (C->D, D->E, ... I->J), then populate the new
# column C with the Industry header + per-row industry values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C - this shifts C:I -> D:J automatically
# (values, formulas and formatting move with it).
$ws.Columns.Item(3).EntireColumn.Insert()

# Header
$ws.Cells.Item(1, 3).Value = "Industry"

$industries = @(
    "Banks",
    "Banks",
    "Metals & Minerals Trading",
    "Finance",
    "Pharmaceuticals & Biotechnology",
    "Power",
    "Insurance",
    "Automobiles",
    "Banks",
    "Beverages",
    "Consumable Fuels",
    "Financial Technology (Fintech)",
    "IT - Software",
    "Capital Markets",
    "Personal Products",
    "Auto Components",
    "Pharmaceuticals & Biotechnology",
    "Chemicals & Petrochemicals",
    "IT - Software",
    "Entertainment",
    "Banks",
    "Petroleum Products",
    "Gas",
    "Consumer Durables",
    "Realty",
    "Power",
    "Pharmaceuticals & Biotechnology",
    "Banks",
    "IT - Services",
    "Agricultural Food & other Products",
    "Diversified FMCG"
)

$row = 2
foreach ($industry in $industries) {
    $ws.Cells.Item($row, 3).Value = $industry
    $row = $row + 1
}
